$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "Create a roadmap, preferrably in a shape form like tree, leaf, Myshow4all Logo, etc., at the top through which any component can be accessed with a referrenced link to each component"

$cell = $ws.Range("A9")
$cell.Value = $newText

# Match the "wrap text" formatting used by the other long story cells (A2, A6)
# so the new row gets the same cell style (wrapText alignment).
$cell.WrapText = $true

# The row needs to grow to fit the wrapped text, same as rows 2 and 6.
$ws.Rows.Item(9).RowHeight = 76.5

$ws.Range("A9").Select()
